$wb = $excel.ActiveWorkbook

# "yh_indexTrend_estimates" was missing a required "date" field. Insert a
# new column B ("date"), shifting growth/period one column to the right:
#   A: index_symbol, B: date, C: growth, D: period
$ws = $wb.Worksheets.Item("yh_indexTrend_estimates")
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "date"

# Copy the existing literal-text "2019-01-01" value (as already used on the
# yh_assetProfile / yh_ohlcv_1d sheets) into the new date column so it stays
# plain text instead of being reinterpreted as a date serial number.
$dateCell = $wb.Worksheets.Item("yh_assetProfile").Range("B2")
$dateCell.Copy($ws.Range("B2"))
$dateCell.Copy($ws.Range("B3"))
$dateCell.Copy($ws.Range("B4"))
$dateCell.Copy($ws.Range("B5"))
